# Generate Report for Handoff
# Updates the localization-status report: swaps the handed-off source file's
# GUID/hash identifiers for the new handoff round, refreshes timestamps, and
# clears the "latest handback" columns (I, J, K) since a new handoff round
# has just started (no handback has happened yet for it).

$wb = $excel.ActiveWorkbook

$oldGuid = "5618eef6-2572-4309-abff-b8645fe5ce31"
$newGuid = "5f6dea56-bf07-4831-acf2-6ec26b6a5adc"
$oldHash = "eedaa5a4e5c001da6d6e901a393d35f43e221077"
$newHash = "3e54f58976117d860e8df3de4230e7ed4b931c05"

$oldMdName  = "$oldGuid.md"
$newMdName  = "$newGuid.md"
$oldMdPath  = "e2e\$oldGuid.md"
$newMdPath  = "e2e\$newGuid.md"

$defaultHandbackDate = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("B2").Value = $newMdPath
$wsOverview.Range("G2").Value = "2016-08-15 16:56:33"

# Re-create the hyperlink on B2 with the new display text, keeping the
# original (external) link target untouched.
$ovB2Address = $null
foreach ($hl in $wsOverview.Hyperlinks) {
    if ($hl.Range.Address() -eq "`$B`$2") {
        $ovB2Address = $hl.Address
    }
}
if (-not $ovB2Address) {
    $ovB2Address = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/501147d9e309f4b6c7e21dee3849dd8e7fd4d96f/e2e/$oldGuid.md"
}
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $ovB2Address, $null, $null, $newMdPath) | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Capture existing hyperlink target addresses before we touch anything.
$zhA2Address = $null
foreach ($hl in $wsZh.Hyperlinks) {
    if ($hl.Range.Address() -eq "`$A`$2") {
        $zhA2Address = $hl.Address
    }
}
if (-not $zhA2Address) {
    $zhA2Address = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/501147d9e309f4b6c7e21dee3849dd8e7fd4d96f/e2e/$oldGuid.md"
}

$wsZh.Range("A2").Value = $newMdName
$wsZh.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-15 16:56:28"
$wsZh.Range("I2").ClearContents()
$wsZh.Range("J2").ClearContents()
$wsZh.Range("K2").Value = $defaultHandbackDate

# Drop every hyperlink on the sheet (A2 + I2) then restore only the one that
# should remain (A2), pointing at the unchanged external URL but with the
# refreshed display text. I2 no longer holds a value, so it no longer gets a
# hyperlink.
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $zhA2Address, $null, $null, $newMdName) | Out-Null

# I2 should no longer look like a hyperlink.
$wsZh.Range("I2").Style = "Normal"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deA2Address = $null
foreach ($hl in $wsDe.Hyperlinks) {
    if ($hl.Range.Address() -eq "`$A`$2") {
        $deA2Address = $hl.Address
    }
}
if (-not $deA2Address) {
    $deA2Address = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/501147d9e309f4b6c7e21dee3849dd8e7fd4d96f/e2e/$oldGuid.md"
}

$wsDe.Range("A2").Value = $newMdName
$wsDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-15 16:56:33"
$wsDe.Range("I2").ClearContents()
$wsDe.Range("J2").ClearContents()
$wsDe.Range("K2").Value = $defaultHandbackDate

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $deA2Address, $null, $null, $newMdName) | Out-Null

$wsDe.Range("I2").Style = "Normal"

# ---------------------------------------------------------------------
# Column width touch-up on zh-cn / de-de sheets (columns I & J shrank once
# their long filename/hyperlink contents were cleared out).
# ---------------------------------------------------------------------
foreach ($ws in @($wsZh, $wsDe)) {
    $ws.Columns.Item(9).ColumnWidth = 18.6506053379604
    $ws.Columns.Item(10).ColumnWidth = 21.7054770333426
}
